$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.487.73'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '2.159.13'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.01'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +3.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.67'
$ws.Range('E7').Value = '  +2.06%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.394'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0851'
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.02'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').Value = '2.480.71'
$ws.Range('E13').Value = '  +2.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.02'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.810'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.51'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '2.164.07'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '39.527.65'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.30'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '0.0₃0846'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.63'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.63'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.17'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.83'
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('E31').Value = '  +5.42%  '
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.59'
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0619'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.63'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.69'
$ws.Range('E40').Value = '  +13.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.46'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '1.528.15'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.19'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0923'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('E47').Value = '  +3.30%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.70'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('E51').Value = '  +35.04%  '
